$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column R: pull in formatting from the matching column-Q cells, then set values/format ---

# R3: same (bottom-border-only) formatting as O3/P3/Q3, left empty.
$ws.Range("Q3").Copy()
$ws.Range("R3").PasteSpecial(-4122)

# R4: 2021 header, same formatting as the other year headers (Q4).
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)
$ws.Range("R4").Value = 2021

# R5:R7 data rows share the thin (no-border) numeric style used by D:Q in those rows.
$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122)
$ws.Range("R5").Value = 20.461775421944708

$ws.Range("Q6").Copy()
$ws.Range("R6").PasteSpecial(-4122)
$ws.Range("R6").Value = 10.464183282203864

$ws.Range("Q7").Copy()
$ws.Range("R7").PasteSpecial(-4122)
$ws.Range("R7").Value = 21.69437772849707

# R8 (total row) uses the bottom-bordered numeric style used by D:Q in row 8.
$ws.Range("Q8").Copy()
$ws.Range("R8").PasteSpecial(-4122)
$ws.Range("R8").Value = 206.4

# --- Re-format the existing yearly figures (rows 5-8, cols D:Q) from "0.00" to "0.0" ---
$ws.Range("D5:R8").NumberFormat = "0.0"

# --- Restore selection shown in the saved workbook ---
$ws.Range("J15").Select()
